# Apply repulled data updates to column F (dSF) for specific rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 1
    6  = -5
    8  = 3
    10 = -3
    12 = -5
    13 = 2
    15 = -7
    21 = -3
    29 = -3
    30 = 5
    31 = -8
    32 = -3
    33 = 2
    36 = -4
    37 = -2
    38 = -5
    41 = 3
    43 = 2
    45 = 3
    46 = 0
    47 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
